# Tilpasset modell for inntak av excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Input"
$ws.Name = "Input"

# Header row
$ws.Range("A1").Value = "PV1"
$ws.Range("B1").Value = "PV2"
$ws.Range("C1").Value = "PV3"
$ws.Range("D1").Value = "Battery1"
$ws.Range("E1").Value = "Battery2"
$ws.Range("F1").Value = "Battery3"

# Row 2
$ws.Range("A2").Value = 100
$ws.Range("B2").Value = 200
$ws.Range("C2").Value = 300
$ws.Range("D2").Value = 500
$ws.Range("E2").Value = 1000
$ws.Range("F2").Value = 2000

# Row 3
$ws.Range("A3").Value = 150
$ws.Range("B3").Value = 250
$ws.Range("C3").Value = 350
$ws.Range("D3").Value = 7000
$ws.Range("E3").Value = 1500
$ws.Range("F3").Value = 3000

# Row 4
$ws.Range("A4").Value = 200
$ws.Range("B4").Value = 300
$ws.Range("C4").Value = 400
$ws.Range("D4").Value = 900
$ws.Range("E4").Value = 2000
$ws.Range("F4").Value = 4000

# Match the author's final selection state
$ws.Range("G9:G10").Select()
